$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-12-12T07:01:17.651136+00:00"
$ws.Range("K3").Value = "2025-12-12T07:01:17.651170+00:00"
$ws.Range("K4").Value = "2025-12-12T07:01:17.651190+00:00"
$ws.Range("K5").Value = "2025-12-12T07:01:19.988360+00:00"
$ws.Range("K6").Value = "2025-12-12T07:01:19.988389+00:00"
$ws.Range("K7").Value = "2025-12-12T07:01:19.988409+00:00"
$ws.Range("K8").Value = "2025-12-12T07:01:22.815018+00:00"
$ws.Range("K9").Value = "2025-12-12T07:01:25.697177+00:00"
$ws.Range("K10").Value = "2025-12-12T07:01:28.681456+00:00"
$ws.Range("K11").Value = "2025-12-12T07:01:31.460695+00:00"
$ws.Range("K12").Value = "2025-12-12T07:01:36.848681+00:00"
$ws.Range("K13").Value = "2025-12-12T07:01:36.848711+00:00"
$ws.Range("K14").Value = "2025-12-12T07:01:39.521990+00:00"
$ws.Range("K15").Value = "2025-12-12T07:01:42.272750+00:00"
$ws.Range("K16").Value = "2025-12-12T07:01:45.043153+00:00"
$ws.Range("K17").Value = "2025-12-12T07:01:47.930946+00:00"
$ws.Range("K18").Value = "2025-12-12T07:01:47.930975+00:00"
$ws.Range("K19").Value = "2025-12-12T07:01:50.681576+00:00"
$ws.Range("K20").Value = "2025-12-12T07:01:50.681604+00:00"
$ws.Range("K21").Value = "2025-12-12T07:01:50.681622+00:00"
$ws.Range("K22").Value = "2025-12-12T07:01:53.443452+00:00"
$ws.Range("K23").Value = "2025-12-12T07:01:53.443480+00:00"
$ws.Range("K24").Value = "2025-12-12T07:01:53.443497+00:00"
$ws.Range("K25").Value = "2025-12-12T07:01:53.443513+00:00"
$ws.Range("K26").Value = "2025-12-12T07:01:53.443529+00:00"
$ws.Range("K27").Value = "2025-12-12T07:01:56.129701+00:00"
$ws.Range("K28").Value = "2025-12-12T07:01:56.129735+00:00"
$ws.Range("K29").Value = "2025-12-12T07:01:56.129754+00:00"
$ws.Range("K30").Value = "2025-12-12T07:01:58.478195+00:00"
$ws.Range("K31").Value = "2025-12-12T07:01:58.478224+00:00"
$ws.Range("K32").Value = "2025-12-12T07:02:00.796802+00:00"
$ws.Range("K33").Value = "2025-12-12T07:02:03.096028+00:00"
$ws.Range("K34").Value = "2025-12-12T07:02:03.096057+00:00"
$ws.Range("K35").Value = "2025-12-12T07:02:03.096076+00:00"
$ws.Range("K36").Value = "2025-12-12T07:02:05.844543+00:00"
$ws.Range("K37").Value = "2025-12-12T07:02:05.844572+00:00"
$ws.Range("K38").Value = "2025-12-12T07:02:08.253503+00:00"
$ws.Range("K39").Value = "2025-12-12T07:02:08.253535+00:00"
$ws.Range("K40").Value = "2025-12-12T07:02:10.946054+00:00"
$ws.Range("K41").Value = "2025-12-12T07:02:10.946081+00:00"
